# Applies targeted numeric updates to the Leve profit-tracking sheets
# (values refreshed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 267629.06
$ws.Range("I98").Value = 350458.03
$ws.Range("J98").Value = 2576.4
$ws.Range("K98").Value = 350458.03
$ws.Range("L98").Value = 2576.4
$ws.Range("M98").Value = -348960.03
$ws.Range("N98").Value = -5572.4
$ws.Range("H122").Value = 267629.06
$ws.Range("I122").Value = 350458.03
$ws.Range("J122").Value = 2576.4
$ws.Range("K122").Value = 1051374.09
$ws.Range("L122").Value = 7729.200000000001
$ws.Range("M122").Value = -1048924.09
$ws.Range("N122").Value = -12629.2
$ws.Range("H133").Value = 29653.545
$ws.Range("J133").Value = 29653.545
$ws.Range("L133").Value = 29653.545
$ws.Range("N133").Value = -39773.545
$ws.Range("H136").Value = 57260
$ws.Range("J136").Value = 57260
$ws.Range("L136").Value = 57260
$ws.Range("N136").Value = -67460
$ws.Range("H137").Value = 41667868
$ws.Range("I137").Value = 50001120
$ws.Range("J137").Value = 1600.5
$ws.Range("K137").Value = 150003360
$ws.Range("L137").Value = 4801.5
$ws.Range("M137").Value = -150000810
$ws.Range("N137").Value = -9901.5
$ws.Range("H138").Value = 5070829
$ws.Range("I138").Value = 2318561.2
$ws.Range("J138").Value = 5955486.5
$ws.Range("K138").Value = 6955683.600000001
$ws.Range("L138").Value = 17866459.5
$ws.Range("M138").Value = -6950543.600000001
$ws.Range("N138").Value = -17876739.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16487.93
$ws.Range("I32").Value = 2654.0317
$ws.Range("J32").Value = 113325.22
$ws.Range("K32").Value = 2654.0317
$ws.Range("L32").Value = 113325.22
$ws.Range("M32").Value = -2367.0317
$ws.Range("N32").Value = -113899.22
$ws.Range("H45").Value = 861.875
$ws.Range("I45").Value = 675
$ws.Range("J45").Value = 924.1667
$ws.Range("K45").Value = 675
$ws.Range("L45").Value = 924.1667
$ws.Range("M45").Value = -298
$ws.Range("N45").Value = -1678.1667
$ws.Range("H61").Value = 1963.875
$ws.Range("I61").Value = 1667
$ws.Range("J61").Value = 2854.5
$ws.Range("K61").Value = 1667
$ws.Range("L61").Value = 2854.5
$ws.Range("M61").Value = -1455
$ws.Range("N61").Value = -3278.5
$ws.Range("H136").Value = 1963.875
$ws.Range("I136").Value = 1667
$ws.Range("J136").Value = 2854.5
$ws.Range("K136").Value = 5001
$ws.Range("L136").Value = 8563.5
$ws.Range("M136").Value = -2451
$ws.Range("N136").Value = -13663.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 28401
$ws.Range("J92").Value = 28401
$ws.Range("L92").Value = 28401
$ws.Range("N92").Value = -33393
$ws.Range("H102").Value = 9556
$ws.Range("I102").Value = 9556
$ws.Range("K102").Value = 9556
$ws.Range("M102").Value = -6311
$ws.Range("H134").Value = 43482396
$ws.Range("I134").Value = 90911416
$ws.Range("J134").Value = 5793.6665
$ws.Range("K134").Value = 272734248
$ws.Range("L134").Value = 17380.9995
$ws.Range("M134").Value = -272731713
$ws.Range("N134").Value = -22450.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 11923.75
$ws.Range("J8").Value = 13570
$ws.Range("L8").Value = 13570
$ws.Range("N8").Value = -13850
$ws.Range("H31").Value = 1298.1025
$ws.Range("I31").Value = 927.2143
$ws.Range("J31").Value = 2242.182
$ws.Range("K31").Value = 927.2143
$ws.Range("L31").Value = 2242.182
$ws.Range("M31").Value = -632.2143
$ws.Range("N31").Value = -2832.182
$ws.Range("H34").Value = 1298.1025
$ws.Range("I34").Value = 927.2143
$ws.Range("J34").Value = 2242.182
$ws.Range("K34").Value = 927.2143
$ws.Range("L34").Value = 2242.182
$ws.Range("M34").Value = -725.2143
$ws.Range("N34").Value = -2646.182
$ws.Range("H93").Value = 13383.25
$ws.Range("I93").Value = 11690.818
$ws.Range("J93").Value = 32000
$ws.Range("K93").Value = 11690.818
$ws.Range("L93").Value = 32000
$ws.Range("M93").Value = -9818.817999999999
$ws.Range("N93").Value = -35744
$ws.Range("H134").Value = 2487.5227
$ws.Range("I134").Value = 1453.1936
$ws.Range("J134").Value = 4954
$ws.Range("K134").Value = 4359.5808
$ws.Range("L134").Value = 14862
$ws.Range("M134").Value = -1824.5808
$ws.Range("N134").Value = -19932

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50000130
$ws.Range("I2").Value = 24.941177
$ws.Range("J2").Value = 333334050
$ws.Range("K2").Value = 149.647062
$ws.Range("L2").Value = 2000004300
$ws.Range("M2").Value = -36.64706200000001
$ws.Range("N2").Value = -2000004526
$ws.Range("H5").Value = 1482.2084
$ws.Range("I5").Value = 625.7778
$ws.Range("J5").Value = 1996.0667
$ws.Range("K5").Value = 1877.3334
$ws.Range("L5").Value = 5988.2001
$ws.Range("M5").Value = -1765.3334
$ws.Range("N5").Value = -6212.2001
$ws.Range("H107").Value = 479.7857
$ws.Range("I107").Value = 506
$ws.Range("J107").Value = 469.3
$ws.Range("K107").Value = 1518
$ws.Range("L107").Value = 1407.9
$ws.Range("M107").Value = 402
$ws.Range("N107").Value = -5247.9
$ws.Range("H135").Value = 1482.2084
$ws.Range("I135").Value = 625.7778
$ws.Range("J135").Value = 1996.0667
$ws.Range("K135").Value = 5632.000199999999
$ws.Range("L135").Value = 17964.6003
$ws.Range("M135").Value = -3097.000199999999
$ws.Range("N135").Value = -23034.6003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1810
$ws.Range("J22").Value = 1810
$ws.Range("L22").Value = 1810
$ws.Range("N22").Value = -2400
$ws.Range("H27").Value = 1810
$ws.Range("J27").Value = 1810
$ws.Range("L27").Value = 1810
$ws.Range("N27").Value = -2024
$ws.Range("H68").Value = 2557.7144
$ws.Range("I68").Value = 2400.8
$ws.Range("J68").Value = 2950
$ws.Range("K68").Value = 2400.8
$ws.Range("L68").Value = 2950
$ws.Range("M68").Value = -1651.8
$ws.Range("N68").Value = -4448
$ws.Range("H71").Value = 2557.7144
$ws.Range("I71").Value = 2400.8
$ws.Range("J71").Value = 2950
$ws.Range("K71").Value = 12004
$ws.Range("L71").Value = 14750
$ws.Range("M71").Value = -8260
$ws.Range("N71").Value = -22238
$ws.Range("H82").Value = 1389.1666
$ws.Range("I82").Value = 1200
$ws.Range("J82").Value = 1427
$ws.Range("K82").Value = 1200
$ws.Range("L82").Value = 1427
$ws.Range("M82").Value = -839
$ws.Range("N82").Value = -2149
$ws.Range("H85").Value = 1389.1666
$ws.Range("I85").Value = 1200
$ws.Range("J85").Value = 1427
$ws.Range("K85").Value = 1200
$ws.Range("L85").Value = 1427
$ws.Range("M85").Value = 48
$ws.Range("N85").Value = -3923

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 29704.2
$ws.Range("I62").Value = 53150
$ws.Range("K62").Value = 53150
$ws.Range("M62").Value = -52526
$ws.Range("H65").Value = 29704.2
$ws.Range("I65").Value = 53150
$ws.Range("K65").Value = 265750
$ws.Range("M65").Value = -262630
